$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aide (Référence)")

# Find the last used row in column A (currently row 23) and add a new
# row right after it with the new reference link, matching the
# formatting ("Lien hypertexte" style) used by the other link rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$url = "https://elifulkerson.com/projects/commandline-wav-player.php"

$cell = $ws.Cells.Item($newRow, 1)
$cell.Value = $url

# Register the hyperlink (this also writes the URL into the cell and
# applies the hyperlink style, but reassigning the named style
# afterwards keeps it consistent with the sibling cells).
$ws.Hyperlinks.Add($cell, $url)
$cell.Style = "Lien hypertexte"
